$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$values = @(
    "0+14=14",
    "72+7=79",
    "31+34=65",
    "19+56=75",
    "98-88=10",
    "19+12=31",
    "29-21=8",
    "17+40=57",
    "30+53=83",
    "76+12=88",
    "38-17=21",
    "59+21=80",
    "81-72=9",
    "6+70=76",
    "91-89=2",
    "14+20=34",
    "82-24=58",
    "90-41=49",
    "39+48=87",
    "77-6=71",
    "44-4=40",
    "14+67=81",
    "98-51=47",
    "57+25=82",
    "24+72=96",
    "92-61=31",
    "59+40=99",
    "2+35=37",
    "58-30=28",
    "31+39=70",
    "51-30=21",
    "57+7=64",
    "4+80=84",
    "41-2=39",
    "64+1=65",
    "40+34=74",
    "58-57=1",
    "84-61=23",
    "9+29=38",
    "67-67=0",
    "29-13=16",
    "82-52=30",
    "99-71=28",
    "29+46=75",
    "74+9=83",
    "94-59=35",
    "82-16=66",
    "69-33=36",
    "72+7=79",
    "94+3=97",
    "27+1=28",
    "29+52=81",
    "47+36=83",
    "43+14=57",
    "75+2=77",
    "36+15=51",
    "1+97=98",
    "5+66=71",
    "3+38=41",
    "95-51=44",
    "29-26=3",
    "49-42=7",
    "26+66=92",
    "43+11=54",
    "75-20=55",
    "73-62=11",
    "27+54=81",
    "74-6=68",
    "14+40=54",
    "62-5=57",
    "33+56=89",
    "54-51=3",
    "60-34=26",
    "62-41=21",
    "61-15=46",
    "33+14=47",
    "97-55=42",
    "3+65=68",
    "48-7=41",
    "38+3=41",
    "47+21=68",
    "44+16=60",
    "30-3=27",
    "26+56=82",
    "91-78=13",
    "19+10=29",
    "37+51=88",
    "73-3=70",
    "51+21=72",
    "8+71=79",
    "41+16=57",
    "46-34=12",
    "1+46=47",
    "26+6=32",
    "65-19=46",
    "17+3=20",
    "25-3=22",
    "43+3=46",
    "85-29=56",
    "92-59=33"
)
$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}
Write-Host "Done: $idx cells updated"
